# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, applied identically on both sheets.
$updates = @{
    2  = 165
    3  = 7382
    4  = 7263
    5  = 90
    11 = 123
    13 = 86
    14 = 659
    15 = 470
    16 = 55
    17 = 21
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
